$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 274
$ws.Range("I4").Value = 93.333336
$ws.Range("J4").Value = 545
$ws.Range("K4").Value = 93.333336
$ws.Range("L4").Value = 545
$ws.Range("M4").Value = 20.666664
$ws.Range("N4").Value = -773

# Row 129
$ws.Range("H129").Value = 955.625
$ws.Range("J129").Value = 1004.1964
$ws.Range("L129").Value = 3012.5892
$ws.Range("N129").Value = -13012.5892

# Row 138
$ws.Range("H138").Value = 3988.6619
$ws.Range("I138").Value = 2028.4186
$ws.Range("J138").Value = 7360.28
$ws.Range("K138").Value = 6085.2558
$ws.Range("L138").Value = 22080.84
$ws.Range("M138").Value = -945.2557999999999
$ws.Range("N138").Value = -32360.84


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 19000
$ws.Range("I37").Value = 19000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 19000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -18727
$ws.Range("N37").ClearContents()

# Row 74
$ws.Range("H74").Value = 1318.6154
$ws.Range("I74").Value = 1287.4546
$ws.Range("J74").Value = 1490
$ws.Range("K74").Value = 1287.4546
$ws.Range("L74").Value = 1490
$ws.Range("M74").Value = -413.4546
$ws.Range("N74").Value = -3238

# Row 77
$ws.Range("H77").Value = 1318.6154
$ws.Range("I77").Value = 1287.4546
$ws.Range("J77").Value = 1490
$ws.Range("K77").Value = 6437.273
$ws.Range("L77").Value = 7450
$ws.Range("M77").Value = -2069.273
$ws.Range("N77").Value = -16186

# Row 122
$ws.Range("H122").Value = 1512.1
$ws.Range("I122").Value = 1263.25
$ws.Range("J122").Value = 2507.5
$ws.Range("K122").Value = 3789.75
$ws.Range("L122").Value = 7522.5
$ws.Range("M122").Value = -1339.75
$ws.Range("N122").Value = -12422.5

# Row 132
$ws.Range("H132").Value = 3085.4783
$ws.Range("I132").Value = 2386
$ws.Range("K132").Value = 7158
$ws.Range("M132").Value = -4628


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2367.1875
$ws.Range("I134").Value = 2060.5518
$ws.Range("J134").Value = 5331.3335
$ws.Range("K134").Value = 6181.655400000001
$ws.Range("L134").Value = 15994.0005
$ws.Range("M134").Value = -3646.655400000001
$ws.Range("N134").Value = -21064.0005


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1997.175
$ws.Range("I31").Value = 1392.3793
$ws.Range("J31").Value = 3591.6365
$ws.Range("K31").Value = 1392.3793
$ws.Range("L31").Value = 3591.6365
$ws.Range("M31").Value = -1097.3793
$ws.Range("N31").Value = -4181.636500000001

# Row 34
$ws.Range("H34").Value = 1997.175
$ws.Range("I34").Value = 1392.3793
$ws.Range("J34").Value = 3591.6365
$ws.Range("K34").Value = 1392.3793
$ws.Range("L34").Value = 3591.6365
$ws.Range("M34").Value = -1190.3793
$ws.Range("N34").Value = -3995.6365

# Row 132
$ws.Range("H132").Value = 323368.5
$ws.Range("I132").Value = 437034.34
$ws.Range("J132").Value = 3037.4546
$ws.Range("K132").Value = 1311103.02
$ws.Range("L132").Value = 9112.363799999999
$ws.Range("M132").Value = -1308573.02
$ws.Range("N132").Value = -14172.3638


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 1428.75
$ws.Range("I2").Value = 2662.5
$ws.Range("J2").Value = 195
$ws.Range("K2").Value = 15975
$ws.Range("L2").Value = 1170
$ws.Range("M2").Value = -15862
$ws.Range("N2").Value = -1396

# Row 113
$ws.Range("H113").Value = 632.4
$ws.Range("I113").Value = 550.375
$ws.Range("J113").Value = 726.1429000000001
$ws.Range("K113").Value = 1651.125
$ws.Range("L113").Value = 2178.4287
$ws.Range("M113").Value = 518.875
$ws.Range("N113").Value = -6518.4287

# Row 122
$ws.Range("H122").Value = 864.2381
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 1088.3846
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 9795.4614
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -14695.4614

# Row 129
$ws.Range("H129").Value = 2274203
$ws.Range("I129").Value = 878.2857
$ws.Range("J129").Value = 3335087.8
$ws.Range("K129").Value = 2634.8571
$ws.Range("L129").Value = 10005263.4
$ws.Range("M129").Value = 2365.1429
$ws.Range("N129").Value = -10015263.4

# Row 131
$ws.Range("H131").Value = 13172672
$ws.Range("I131").Value = 12462.333
$ws.Range("J131").Value = 14940462
$ws.Range("K131").Value = 37386.999
$ws.Range("L131").Value = 44821386
$ws.Range("M131").Value = -32346.999
$ws.Range("N131").Value = -44831466

# Row 137
$ws.Range("H137").Value = 19610950
$ws.Range("I137").Value = 2261.4285
$ws.Range("J137").Value = 33337030
$ws.Range("K137").Value = 6784.2855
$ws.Range("L137").Value = 100011090
$ws.Range("M137").Value = -1684.2855
$ws.Range("N137").Value = -100021290


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2818.9
$ws.Range("I122").Value = 2173
$ws.Range("J122").Value = 4326
$ws.Range("K122").Value = 6519
$ws.Range("L122").Value = 12978
$ws.Range("M122").Value = -4069
$ws.Range("N122").Value = -17878

# Row 132
$ws.Range("H132").Value = 1188.5405
$ws.Range("I132").Value = 753.1539
$ws.Range("J132").Value = 2217.6365
$ws.Range("K132").Value = 2259.4617
$ws.Range("L132").Value = 6652.9095
$ws.Range("M132").Value = 270.5383000000002
$ws.Range("N132").Value = -11712.9095


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 30814.572
$ws.Range("I61").Value = 51825.75
$ws.Range("K61").Value = 51825.75
$ws.Range("M61").Value = -51623.75

# Row 113
$ws.Range("H113").Value = 30814.572
$ws.Range("I113").Value = 51825.75
$ws.Range("K113").Value = 51825.75
$ws.Range("M113").Value = -49655.75

# Row 132
$ws.Range("H132").Value = 3230.875
$ws.Range("I132").Value = 2741.8333
$ws.Range("J132").Value = 4698
$ws.Range("K132").Value = 8225.499899999999
$ws.Range("L132").Value = 14094
$ws.Range("M132").Value = -5695.499899999999
$ws.Range("N132").Value = -19154


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1358.5834
$ws.Range("I132").Value = 977.8148
$ws.Range("J132").Value = 2500.889
$ws.Range("K132").Value = 2933.4444
$ws.Range("L132").Value = 7502.667
$ws.Range("M132").Value = -403.4443999999999
$ws.Range("N132").Value = -12562.667

# Row 136
$ws.Range("H136").Value = 1704.0588
$ws.Range("I136").Value = 1494.6
$ws.Range("J136").Value = 3275
$ws.Range("K136").Value = 4483.799999999999
$ws.Range("L136").Value = 9825
$ws.Range("M136").Value = -1933.799999999999
$ws.Range("N136").Value = -14925

